# Fig-5-Supp-Table-2-data.xlsx edit script
# 1) Rename 4 existing sheets ("Fig 5 ..." -> "Fig 5A ...")
# 2) Insert a brand-new worksheet "Fig 5C ZIKV IgM ELISA " between the Fig 5A
#    sheets and the Supp Table 2 sheets, and populate it with the ELISA data.
# 3) Make the new sheet the active / selected sheet (tabSelected + activeTab).
# 4) Turn on iterative calculation for the workbook.

$wb = $excel.ActiveWorkbook

# --- 1) Rename the existing sheets -----------------------------------------
$wb.Worksheets.Item("Fig 5 ZIKV-PR(044) MFI raw data").Name  = "Fig 5A ZIKV-PR(044) raw data"
$wb.Worksheets.Item("Fig 5 ZIKV-PR MFI figure data ").Name   = "Fig 5A ZIKV-PR MFI figure data "
$wb.Worksheets.Item("Fig5 ZIKV-DAK(030) MFI raw data").Name  = "Fig 5A ZIKV-DAK(030) raw data"
$wb.Worksheets.Item("Fig 5 ZIKV-DAK MFI figure data").Name   = "Fig 5A ZIKV-DAK MFI figure data"

# --- 2) Insert the new ELISA worksheet --------------------------------------
$afterSheet = $wb.Worksheets.Item("Fig 5A ZIKV-DAK MFI figure data")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "Fig 5C ZIKV IgM ELISA "

# Header row
$newSheet.Range("A1").Value = "Sample"
$newSheet.Range("B1").Value = "Collection Time"
$newSheet.Range("C1").Value = "Replicate 1 Extinction Value"
$newSheet.Range("D1").Value = "Replicate 2 Extinction Value "
$newSheet.Range("E1").Value = "Average"

# Data rows
$newSheet.Range("A2").Value = "030-501"
$newSheet.Range("B2").Value = "Day of Delivery"
$newSheet.Range("C2").Value = 0.25
$newSheet.Range("D2").Value = 0.28
$newSheet.Range("E2").Value = 0.265

$newSheet.Range("A3").Value = "030-502"
$newSheet.Range("B3").Value = "Day of Delivery"
$newSheet.Range("C3").Value = 0.27
$newSheet.Range("D3").Value = 0.27
$newSheet.Range("E3").Value = 0.27

$newSheet.Range("A4").Value = "030-503"
$newSheet.Range("B4").Value = "Day of Delivery"
$newSheet.Range("C4").Value = 0.27
$newSheet.Range("D4").Value = 0.24
$newSheet.Range("E4").Value = 0.255

$newSheet.Range("A5").Value = "030-504"
$newSheet.Range("B5").Value = "Day of Delivery"
$newSheet.Range("C5").Value = 0.29
$newSheet.Range("D5").Value = 0.25
$newSheet.Range("E5").Value = 0.27

$newSheet.Range("A6").Value = "030-101"
$newSheet.Range("B6").Value = "4 dpi"
$newSheet.Range("C6").Value = 0.5
$newSheet.Range("D6").Value = 0.55
$newSheet.Range("E6").Value = 0.525

$newSheet.Range("A7").Value = "030-101"
$newSheet.Range("B7").Value = "14 dpi"
$newSheet.Range("C7").Value = 10.46
$newSheet.Range("D7").Value = 8.99
$newSheet.Range("E7").Value = 9.725

$newSheet.Range("A8").Value = "Assay Postive "
$newSheet.Range("B8").Value = "n/a"
$newSheet.Range("C8").Value = 3.46

$newSheet.Range("A9").Value = "Assay Negative Control"
$newSheet.Range("B9").Value = "n/a"
$newSheet.Range("C9").Value = 0.46

# Match the font used by the rest of the new sheet's cells (size 12, black)
$newSheet.Range("A1:E9").Font.Size = 12
$newSheet.Range("A1:E9").Font.Color = 0
$newSheet.Range("A1:E9").RowHeight = 16

# --- 3) Make the new sheet the active tab -----------------------------------
$newSheet.Activate()
[void]$newSheet.Range("O38").Select()

# --- 4) Iterative calculation -----------------------------------------------
$excel.Iteration = $true
